# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.154.25'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '3.604.04'
$ws.Range('E3').Value = '  +2.59%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''604.57'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = '''196.73'
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('D7').Value = '''0.625'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.206'
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').Value = '''0.649'
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '''0.0000304'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = '''9.57'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '4.177.27'
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('D15').Value = '''13.09'
$ws.Range('E15').Value = '  +4.04%  '
$ws.Range('D16').Value = '''594.40'
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').Value = '70.285.72'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').Value = '''19.12'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '3.592.73'
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('E20').Value = '  +1.49%  '
$ws.Range('D21').Value = '''0.995'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = '''17.71'
$ws.Range('E22').Value = '  -2.20%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '''101.54'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('E26').Value = '  -1.54%  '
$ws.Range('D27').Value = '''10.74'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').Value = '''33.84'
$ws.Range('E29').Value = '  +0.94%  '
$ws.Range('D30').Value = '''4.72'
$ws.Range('E30').Value = '  +3.76%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').Value = '''12.29'
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').Value = '''63.25'
$ws.Range('D35').Value = '0.0₃0886'
$ws.Range('E35').Value = '  +8.80%  '
$ws.Range('D36').Value = '3.942.81'
$ws.Range('E36').Value = '  +5.44%  '
$ws.Range('D37').Value = '''3.11'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''524.53'
$ws.Range('E38').Value = '  +7.89%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '''1.00'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('D40').Value = '''36.84'
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').Value = '''0.390'
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').Value = '''0.133'
$ws.Range('E43').Value = '  -1.99%  '
$ws.Range('D44').Value = '''0.0454'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '''3.44'
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('D46').Value = '''2.85'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').Value = '''0.000250'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('E51').Value = '  +2.69%  '
